$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.9

# Row 3
$ws.Range("M3").Value = 1.06
$ws.Range("T3").Value = 1.86
$ws.Range("U3").Value = 1.7

# Row 4
$ws.Range("F4").Value = 1.99
$ws.Range("J4").Value = 3.75
$ws.Range("K4").Value = 3.8
$ws.Range("L4").Value = 1.36
$ws.Range("N4").Value = 4
$ws.Range("O4").Value = 1.31
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1.96
$ws.Range("S4").Value = 3.4
$ws.Range("V4").Value = 1.3
$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 15
$ws.Range("AD4").Value = 16.5
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 10
$ws.Range("AH4").Value = 18.5
$ws.Range("AK4").Value = 20
$ws.Range("AM4").Value = 110
$ws.Range("AO4").Value = 55

# Row 5
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 13.5
$ws.Range("J5").Value = 5.6
$ws.Range("L5").Value = 1.38
$ws.Range("P5").Value = 2
$ws.Range("V5").Value = 1.08
$ws.Range("W5").Value = 3.8
$ws.Range("X5").Value = 16
$ws.Range("Z5").Value = 120
$ws.Range("AA5").Value = 740
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 48
$ws.Range("AE5").Value = 300
$ws.Range("AF5").Value = 6.8
$ws.Range("AH5").Value = 38
$ws.Range("AI5").Value = 240
$ws.Range("AJ5").Value = 9.800000000000001
$ws.Range("AK5").Value = 16
$ws.Range("AL5").Value = 55
$ws.Range("AM5").Value = 290
$ws.Range("AO5").Value = 500

# Row 6
$ws.Range("H6").Value = 2.52
$ws.Range("I6").Value = 2.82
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 1.98
$ws.Range("O6").Value = 1.3
$ws.Range("R6").Value = 1.32
$ws.Range("S6").Value = 2.74
$ws.Range("T6").Value = 1.56
$ws.Range("U6").Value = 1.01
$ws.Range("V6").Value = 1.54
$ws.Range("W6").Value = 1.46
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 16.5
$ws.Range("Z6").Value = 25
$ws.Range("AA6").Value = 1000
$ws.Range("AB6").Value = 18
$ws.Range("AC6").Value = 11
$ws.Range("AD6").Value = 17
$ws.Range("AE6").Value = 38
$ws.Range("AF6").Value = 29
$ws.Range("AG6").Value = 18.5
$ws.Range("AH6").Value = 23
$ws.Range("AI6").Value = 50
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 44
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 1000

# Row 7
$ws.Range("I7").Value = 7.6
$ws.Range("L7").Value = 1.01
$ws.Range("S7").Value = 2.84
$ws.Range("V7").Value = 1.15
$ws.Range("W7").Value = 2.4

# Row 8
$ws.Range("I8").Value = 1.28
$ws.Range("J8").Value = 7.6
$ws.Range("T8").Value = 1.64
$ws.Range("AM8").Value = 80

# Row 9
$ws.Range("I9").Value = 4.3
$ws.Range("Q9").Value = 1.63
$ws.Range("S9").Value = 2.58
$ws.Range("AC9").Value = 10

# Row 10
$ws.Range("G10").Value = 1.81
$ws.Range("H10").Value = 4.6

# Row 11
$ws.Range("H11").Value = 27
$ws.Range("K11").Value = 10.5
$ws.Range("R11").Value = 1.66
$ws.Range("S11").Value = 2.46
$ws.Range("T11").Value = 2.82
$ws.Range("Z11").Value = 410
$ws.Range("AC11").Value = 23
$ws.Range("AG11").Value = 14.5
$ws.Range("AH11").Value = 65

# Row 12
$ws.Range("F12").Value = 1.45
$ws.Range("G12").Value = 1.47
$ws.Range("J12").Value = 4.7
$ws.Range("K12").Value = 4.9

# Row 13
$ws.Range("G13").Value = 7.6
$ws.Range("H13").Value = 1.5
$ws.Range("I13").Value = 1.53

# Row 14
$ws.Range("F14").Value = 1.31
$ws.Range("H14").Value = 7.4
$ws.Range("K14").Value = 10
$ws.Range("N14").Value = 2.16
$ws.Range("P14").Value = 2.16
$ws.Range("Q14").Value = 1.56
$ws.Range("R14").Value = 1.52
$ws.Range("S14").Value = 2.32
$ws.Range("T14").Value = 1.9
$ws.Range("U14").Value = 1.94
